# Automatische test-sync: 2025-06-22 18:55:50
#
# A new inbound mail ("Aanmelding nieuwsbrief") was logged at 18:55:11 and
# classified as "Afmelding / Nieuwsbrief". This appends that entry to the
# "Logs" sheet, extends the two conditional-formatting ranges on that sheet
# to cover the new row, and refreshes the category totals/order on the
# "Dashboard" sheet to reflect the updated counts.

$wb = $excel.ActiveWorkbook

$logs      = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new mail-log entry to the "Logs" sheet (row 25) ---
$newRow = 25
$logs.Cells.Item($newRow, 1).Value = "Aanmelding nieuwsbrief"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik wil me graag inschrijven voor de nieuwsbrief."
$logs.Cells.Item($newRow, 4).Value = "Afmelding / Nieuwsbrief"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 18:55:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$categoryRules = $logs.Range("D2:D24").FormatConditions
for ($i = 1; $i -le $categoryRules.Count; $i++) {
    $categoryRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D25"))
}

$answeredRules = $logs.Range("G2:G24").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G25"))
}

# --- Refresh the "Dashboard" category breakdown (re-sorted by count desc) ---
$dashboard.Cells.Item(4, 1).Value = "Afmelding / Nieuwsbrief"
$dashboard.Cells.Item(4, 2).Value = 3
$dashboard.Cells.Item(5, 1).Value = "Sollicitatie / Vacature"
$dashboard.Cells.Item(5, 2).Value = 2
$dashboard.Cells.Item(6, 1).Value = "Samenwerking / Partnerverzoek"
$dashboard.Cells.Item(6, 2).Value = 2
$dashboard.Cells.Item(8, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(8, 2).Value = 2
